$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (shifts old B:E -> E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# New header row values: B1=Jun_27, C1=Jun_26, D1=Jun_26
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill B2:D27 with "UN" (unchanged) like the rest of the table
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Row 5 (Zacks Investment Research) got an actual rating change instead of "Unchanged"
$ws.Range("B5").Value = "6/27/2018,Upgrades,Sell -> Hold,"
$ws.Range("B5").Interior.ColorIndex = 42

# Two new analyst rows appended at the bottom, only populated for the new date columns
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
